# Apply the "adding data to graphics" edit:
#  - Type sheet: reshuffle category rows 2-6 (new values + new category assignment)
#  - Type+Visibility sheet: update "public" totals for the four groups and add
#    new "Total Static"/"Total Dynamic" roll-up columns (I/J) summing each group
#  - Update selections / active sheet to match the final authoring state

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Type" (A1:H6) - categories got re-ordered / re-valued
# ---------------------------------------------------------------------------
$wsType = $wb.Worksheets.Item("Type")

$wsType.Range("A2").Value = "Constructor Parameter"
$wsType.Range("C2").Value = 33734
$wsType.Range("D2").Value = 3039

$wsType.Range("A3").Value = "Method Parameter"
$wsType.Range("C3").Value = 57869
$wsType.Range("D3").Value = 12912

$wsType.Range("A4").Value = "Field"
$wsType.Range("C4").Value = 6822
$wsType.Range("D4").Value = 1868

$wsType.Range("A5").Value = "Method Return"
$wsType.Range("C5").Value = 109705
$wsType.Range("D5").Value = 52222

$wsType.Range("A6").Value = "Local Variable"
$wsType.Range("C6").Value = 39168
$wsType.Range("D6").Value = 94960

# ---------------------------------------------------------------------------
# Sheet "Type+Visibility" (A1:H26 -> A1:J26) - new data + roll-up columns
# ---------------------------------------------------------------------------
$wsTV = $wb.Worksheets.Item("Type+Visibility")

# New header labels for the roll-up columns
$wsTV.Range("I1").Value = "Total Static"
$wsTV.Range("J1").Value = "Total Dynamic"

# Group 1 (Method Return): rows 2-4, "public" row is row 2
$wsTV.Range("C2").Value = 103503
$wsTV.Range("D2").Value = 50391
$wsTV.Range("I2").Formula = "=SUM(C2:C4)"
$wsTV.Range("J2").Formula = "=SUM(D2:D4)"

# Group 2 (Method Parameter): rows 6-8, "public" row is row 6
$wsTV.Range("C6").Value = 53435
$wsTV.Range("D6").Value = 11430
$wsTV.Range("I6").Formula = "=SUM(C6:C8)"
$wsTV.Range("J6").Formula = "=SUM(D6:D8)"

# "private" row 7 and "protected" row 8 of group 2 also changed
$wsTV.Range("C7").Value = 2386
$wsTV.Range("D7").Value = 1166
$wsTV.Range("D8").Value = 316

# Group 3 (Constructor Parameter): rows 10-12, "public" row is row 10
$wsTV.Range("C10").Value = 33635
$wsTV.Range("D10").Value = 3006
$wsTV.Range("I10").Formula = "=SUM(C10:C12)"
$wsTV.Range("J10").Formula = "=SUM(D10:D12)"

# Group 4 (Field): rows 14-16, "public" row is row 14
$wsTV.Range("I14").Formula = "=SUM(C14:C16)"
$wsTV.Range("J14").Formula = "=SUM(D14:D16)"

# ---------------------------------------------------------------------------
# View state: active sheet / selections
# ---------------------------------------------------------------------------

# Leave a selection behind on "Type" at D3 (matches final saved state)
$wsType.Activate()
$wsType.Range("D3").Select()

# "Project Size" keeps its own last selection (R49); just visit it so it's
# no longer the active tab when we're done (tabSelected moves off it).
$wsProjSize = $wb.Worksheets.Item("Project Size")
$wsProjSize.Activate()

# Final active sheet/selection: "Type+Visibility", covering the new I10:J14 block
$wsTV.Activate()
$wsTV.Range("I10:J14").Select()
